$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Paul George"
$ws.Range("B7").Value = "SG,SF,PF"
$ws.Range("C7").Value = "Philadelphia 76ers"

$ws.Range("A8").Value = "Wendell Carter Jr."
$ws.Range("B8").Value = "PF,C"
$ws.Range("C8").Value = "Orlando Magic"

$ws.Range("A9").Value = "Bilal Coulibaly"
$ws.Range("B9").Value = "SG,SF"
$ws.Range("C9").Value = "Washington Wizards"

$ws.Range("A14").Value = "Jaren Jackson Jr."
$ws.Range("B14").Value = "PF,C"
$ws.Range("C14").Value = "Memphis Grizzlies"

$ws.Range("A15").Value = "Giannis Antetokounmpo"
$ws.Range("B15").Value = "PF,C"
$ws.Range("C15").Value = "Milwaukee Bucks"

$ws.Range("A16").Value = "Ivica Zubac"
$ws.Range("B16").Value = "C"
$ws.Range("C16").Value = "LA Clippers"
